$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target value is textual (stored as inline string in the original
# workbook) even though it looks numeric/percentage. Setting NumberFormat
# to "@" (Text) before assigning the value prevents Excel from silently
# re-interpreting the string as a number/percentage, and resetting the
# style back to "Normal" afterwards keeps the cell formatting identical
# to the original (General format, default style).

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "313.24"
Set-TextValue "E2" "-0.56%"
Set-TextValue "G2" "5"
Set-TextValue "D3" "37.91"
Set-TextValue "E3" "-3.55%"
Set-TextValue "G3" "5"
Set-TextValue "D4" "5.064"
Set-TextValue "E4" "-1.66%"
Set-TextValue "G4" "5"
Set-TextValue "D5" "0.07764"
Set-TextValue "E5" "-4.94%"
Set-TextValue "G5" "5"
Set-TextValue "D6" "4.342"
Set-TextValue "E6" "-0.90%"
Set-TextValue "G6" "5"
Set-TextValue "D7" "1.900"
Set-TextValue "E7" "-4.82%"
Set-TextValue "G7" "5"
Set-TextValue "D8" "8.184"
Set-TextValue "E8" "-1.94%"
Set-TextValue "G8" "5"
Set-TextValue "D9" "0.9169"
Set-TextValue "E9" "-2.27%"
Set-TextValue "G9" "5"
Set-TextValue "D10" "0.1237"
Set-TextValue "E10" "-5.58%"
Set-TextValue "G10" "5"
Set-TextValue "D11" "0.1890"
Set-TextValue "E11" "-3.98%"
Set-TextValue "G11" "5"
Set-TextValue "D12" "0.08841"
Set-TextValue "E12" "-2.86%"
Set-TextValue "G12" "5"
Set-TextValue "D13" "0.03390"
Set-TextValue "E13" "-3.88%"
Set-TextValue "G13" "5"
Set-TextValue "D14" "0.09700"
Set-TextValue "E14" "-0.54%"
Set-TextValue "G14" "5"
Set-TextValue "D15" "0.001365"
Set-TextValue "E15" "-3.31%"
Set-TextValue "G15" "5"
Set-TextValue "D16" "0.005929"
Set-TextValue "E16" "-11.59%"
Set-TextValue "G16" "5"
Set-TextValue "D17" "3.535"
Set-TextValue "E17" "-2.68%"
Set-TextValue "G17" "5"
Set-TextValue "D18" "3.099"
Set-TextValue "E18" "-0.64%"
Set-TextValue "G18" "5"
Set-TextValue "D19" "0.3409"
Set-TextValue "E19" "-1.78%"
Set-TextValue "G19" "5"
Set-TextValue "D20" "0.1296"
Set-TextValue "E20" "-1.60%"
Set-TextValue "G20" "5"
Set-TextValue "D21" "5.028"
Set-TextValue "E21" "0.88%"
Set-TextValue "G21" "5"
Set-TextValue "E22" "4.03%"
Set-TextValue "G22" "5"
Set-TextValue "E23" "5,587.58%"
Set-TextValue "G23" "5"
Set-TextValue "D24" "0.04389"
Set-TextValue "E24" "0.36%"
Set-TextValue "G24" "5"
Set-TextValue "E25" "-2.40%"
Set-TextValue "G25" "5"
Set-TextValue "D26" "0.004248"
Set-TextValue "E26" "-10.73%"
Set-TextValue "G26" "5"
Set-TextValue "E27" "-65.34%"
Set-TextValue "G27" "5"
Set-TextValue "G28" "5"
Set-TextValue "G29" "5"
Set-TextValue "G30" "5"
Set-TextValue "G31" "5"
Set-TextValue "G32" "5"
Set-TextValue "G33" "5"
Set-TextValue "G34" "5"
Set-TextValue "G35" "5"
Set-TextValue "G36" "5"
Set-TextValue "G37" "5"
Set-TextValue "G38" "5"
Set-TextValue "D39" "0.02145"
Set-TextValue "E39" "-4.29%"
Set-TextValue "G39" "5"
Set-TextValue "D40" "0.04980"
Set-TextValue "G40" "5"
Set-TextValue "D41" "0.007764"
Set-TextValue "E41" "0.18%"
Set-TextValue "G41" "5"
Set-TextValue "D42" "0.009956"
Set-TextValue "E42" "-3.53%"
Set-TextValue "G42" "5"
Set-TextValue "E43" "-3.86%"
Set-TextValue "G43" "5"
Set-TextValue "E44" "-2.04%"
Set-TextValue "G44" "5"
Set-TextValue "E45" "9.05%"
Set-TextValue "G45" "5"
Set-TextValue "D46" "0.00006516"
Set-TextValue "E46" "-4.52%"
Set-TextValue "G46" "5"
Set-TextValue "G47" "5"
Set-TextValue "D48" "0.003199"
Set-TextValue "E48" "6.42%"
Set-TextValue "G48" "5"
Set-TextValue "E49" "-0.19%"
Set-TextValue "G49" "5"
Set-TextValue "G50" "5"
Set-TextValue "G51" "5"
